$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header labels in row 1 ("back"/"front" -> "Question"/"Answer")
$ws.Range("A1").Value = "Question"
$ws.Range("B1").Value = "Answer"

# Move the active selection to D5 (as recorded in the saved sheet view)
$ws.Range("D5").Select()
